$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C5")

# New daily-entry data block, rows 7-14, columns A-L
# (values written in plain decimal form - equal to the scientific-notation
# doubles in the target file, e.g. 4.2000000000000003E-2 == 0.042)
$data = @(
    @(0.042, 0.039, 0.047, 0.04, 0.142, 0.043, 0.04, 0.045, 0.04, 0.04, 0.039, 0.039),
    @(0.05, 0.045, 0.044, 0.048, 0.049, 0.042, 0.047, 0.047, 0.047, 0.048, 0.047, 0.04),
    @(0.042, 0.037, 0.046, 0.044, 0.041, 0.041, 0.047, 0.048, 0.045, 0.049, 0.046, 0.072),
    @(0.04, 0.045, 0.043, 0.044, 0.044, 0.044, 0.049, 0.047, 0.048, 0.048, 0.054, 0.042),
    @(0.04, 0.05, 0.044, 0.044, 0.045, 0.044, 0.06, 0.06, 0.062, 0.066, 0.061, 0.04),
    @(0.039, 0.041, 0.042, 0.042, 0.048, 0.043, 0.054, 0.058, 0.058, 0.06, 0.053, 0.039),
    @(0.039, 0.042, 0.044, 0.046, 0.044, 0.043, 0.047, 0.056, 0.056, 0.066, 0.062, 0.039),
    @(0.04, 0.04, 0.043, 0.041, 0.04, 0.039, 0.041, 0.04, 0.04, 0.041, 0.064, 0.039)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 7 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

$ws.Range("A7:L14").Select()
